# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" table (rows 16-105, columns C:G) is rebuilt:
#   - Previously the 88 monthly periods (1610..2402, skipping the
#     non-existent 1703) for worker VIVIANA RODRIGUEZ CORREA (CC 1143325562)
#     were listed in ascending order, with two of her rows (for periods
#     1611 and 1612) having been overwritten by worker EDGAR ALEXANDER
#     TAUTIVA CASALLAS (CC 80007756).
#   - Now all 88 periods for VIVIANA are listed in descending order
#     (2402 down to 1610) in rows 16-103, and EDGAR's two periods
#     (1612, 1611) are appended at the bottom in rows 104-105 with his
#     own "Valor Mora"/"Salario Basico" figures.
#
# Row 16's period (2402) keeps the special partial "Valor Mora" of 21840
# (rather than the usual 36400) -- the special amount stays attached to
# the 2402 period, not to the row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(16,"1143325562","VIVIANA RODRIGUEZ CORREA","2402",21840,910000),
    @(17,"1143325562","VIVIANA RODRIGUEZ CORREA","2401",36400,910000),
    @(18,"1143325562","VIVIANA RODRIGUEZ CORREA","2312",36400,910000),
    @(19,"1143325562","VIVIANA RODRIGUEZ CORREA","2311",36400,910000),
    @(20,"1143325562","VIVIANA RODRIGUEZ CORREA","2310",36400,910000),
    @(21,"1143325562","VIVIANA RODRIGUEZ CORREA","2309",36400,910000),
    @(22,"1143325562","VIVIANA RODRIGUEZ CORREA","2308",36400,910000),
    @(23,"1143325562","VIVIANA RODRIGUEZ CORREA","2307",36400,910000),
    @(24,"1143325562","VIVIANA RODRIGUEZ CORREA","2306",36400,910000),
    @(25,"1143325562","VIVIANA RODRIGUEZ CORREA","2305",36400,910000),
    @(26,"1143325562","VIVIANA RODRIGUEZ CORREA","2304",36400,910000),
    @(27,"1143325562","VIVIANA RODRIGUEZ CORREA","2303",36400,910000),
    @(28,"1143325562","VIVIANA RODRIGUEZ CORREA","2302",36400,910000),
    @(29,"1143325562","VIVIANA RODRIGUEZ CORREA","2301",36400,910000),
    @(30,"1143325562","VIVIANA RODRIGUEZ CORREA","2212",36400,910000),
    @(31,"1143325562","VIVIANA RODRIGUEZ CORREA","2211",36400,910000),
    @(32,"1143325562","VIVIANA RODRIGUEZ CORREA","2210",36400,910000),
    @(33,"1143325562","VIVIANA RODRIGUEZ CORREA","2209",36400,910000),
    @(34,"1143325562","VIVIANA RODRIGUEZ CORREA","2208",36400,910000),
    @(35,"1143325562","VIVIANA RODRIGUEZ CORREA","2207",36400,910000),
    @(36,"1143325562","VIVIANA RODRIGUEZ CORREA","2206",36400,910000),
    @(37,"1143325562","VIVIANA RODRIGUEZ CORREA","2205",36400,910000),
    @(38,"1143325562","VIVIANA RODRIGUEZ CORREA","2204",36400,910000),
    @(39,"1143325562","VIVIANA RODRIGUEZ CORREA","2203",36400,910000),
    @(40,"1143325562","VIVIANA RODRIGUEZ CORREA","2202",36400,910000),
    @(41,"1143325562","VIVIANA RODRIGUEZ CORREA","2201",36400,910000),
    @(42,"1143325562","VIVIANA RODRIGUEZ CORREA","2112",36400,910000),
    @(43,"1143325562","VIVIANA RODRIGUEZ CORREA","2111",36400,910000),
    @(44,"1143325562","VIVIANA RODRIGUEZ CORREA","2110",36400,910000),
    @(45,"1143325562","VIVIANA RODRIGUEZ CORREA","2109",36400,910000),
    @(46,"1143325562","VIVIANA RODRIGUEZ CORREA","2108",36400,910000),
    @(47,"1143325562","VIVIANA RODRIGUEZ CORREA","2107",36400,910000),
    @(48,"1143325562","VIVIANA RODRIGUEZ CORREA","2106",36400,910000),
    @(49,"1143325562","VIVIANA RODRIGUEZ CORREA","2105",36400,910000),
    @(50,"1143325562","VIVIANA RODRIGUEZ CORREA","2104",36400,910000),
    @(51,"1143325562","VIVIANA RODRIGUEZ CORREA","2103",36400,910000),
    @(52,"1143325562","VIVIANA RODRIGUEZ CORREA","2102",36400,910000),
    @(53,"1143325562","VIVIANA RODRIGUEZ CORREA","2101",36400,910000),
    @(54,"1143325562","VIVIANA RODRIGUEZ CORREA","2012",36400,910000),
    @(55,"1143325562","VIVIANA RODRIGUEZ CORREA","2011",36400,910000),
    @(56,"1143325562","VIVIANA RODRIGUEZ CORREA","2010",36400,910000),
    @(57,"1143325562","VIVIANA RODRIGUEZ CORREA","2009",36400,910000),
    @(58,"1143325562","VIVIANA RODRIGUEZ CORREA","2008",36400,910000),
    @(59,"1143325562","VIVIANA RODRIGUEZ CORREA","2007",36400,910000),
    @(60,"1143325562","VIVIANA RODRIGUEZ CORREA","2006",36400,910000),
    @(61,"1143325562","VIVIANA RODRIGUEZ CORREA","2005",36400,910000),
    @(62,"1143325562","VIVIANA RODRIGUEZ CORREA","2004",36400,910000),
    @(63,"1143325562","VIVIANA RODRIGUEZ CORREA","2003",36400,910000),
    @(64,"1143325562","VIVIANA RODRIGUEZ CORREA","2002",36400,910000),
    @(65,"1143325562","VIVIANA RODRIGUEZ CORREA","2001",36400,910000),
    @(66,"1143325562","VIVIANA RODRIGUEZ CORREA","1912",36400,910000),
    @(67,"1143325562","VIVIANA RODRIGUEZ CORREA","1911",36400,910000),
    @(68,"1143325562","VIVIANA RODRIGUEZ CORREA","1910",36400,910000),
    @(69,"1143325562","VIVIANA RODRIGUEZ CORREA","1909",36400,910000),
    @(70,"1143325562","VIVIANA RODRIGUEZ CORREA","1908",36400,910000),
    @(71,"1143325562","VIVIANA RODRIGUEZ CORREA","1907",36400,910000),
    @(72,"1143325562","VIVIANA RODRIGUEZ CORREA","1906",36400,910000),
    @(73,"1143325562","VIVIANA RODRIGUEZ CORREA","1905",36400,910000),
    @(74,"1143325562","VIVIANA RODRIGUEZ CORREA","1904",36400,910000),
    @(75,"1143325562","VIVIANA RODRIGUEZ CORREA","1903",36400,910000),
    @(76,"1143325562","VIVIANA RODRIGUEZ CORREA","1902",36400,910000),
    @(77,"1143325562","VIVIANA RODRIGUEZ CORREA","1901",36400,910000),
    @(78,"1143325562","VIVIANA RODRIGUEZ CORREA","1812",36400,910000),
    @(79,"1143325562","VIVIANA RODRIGUEZ CORREA","1811",36400,910000),
    @(80,"1143325562","VIVIANA RODRIGUEZ CORREA","1810",36400,910000),
    @(81,"1143325562","VIVIANA RODRIGUEZ CORREA","1809",36400,910000),
    @(82,"1143325562","VIVIANA RODRIGUEZ CORREA","1808",36400,910000),
    @(83,"1143325562","VIVIANA RODRIGUEZ CORREA","1807",36400,910000),
    @(84,"1143325562","VIVIANA RODRIGUEZ CORREA","1806",36400,910000),
    @(85,"1143325562","VIVIANA RODRIGUEZ CORREA","1805",36400,910000),
    @(86,"1143325562","VIVIANA RODRIGUEZ CORREA","1804",36400,910000),
    @(87,"1143325562","VIVIANA RODRIGUEZ CORREA","1803",36400,910000),
    @(88,"1143325562","VIVIANA RODRIGUEZ CORREA","1802",36400,910000),
    @(89,"1143325562","VIVIANA RODRIGUEZ CORREA","1801",36400,910000),
    @(90,"1143325562","VIVIANA RODRIGUEZ CORREA","1712",36400,910000),
    @(91,"1143325562","VIVIANA RODRIGUEZ CORREA","1711",36400,910000),
    @(92,"1143325562","VIVIANA RODRIGUEZ CORREA","1710",36400,910000),
    @(93,"1143325562","VIVIANA RODRIGUEZ CORREA","1709",36400,910000),
    @(94,"1143325562","VIVIANA RODRIGUEZ CORREA","1708",36400,910000),
    @(95,"1143325562","VIVIANA RODRIGUEZ CORREA","1707",36400,910000),
    @(96,"1143325562","VIVIANA RODRIGUEZ CORREA","1706",36400,910000),
    @(97,"1143325562","VIVIANA RODRIGUEZ CORREA","1705",36400,910000),
    @(98,"1143325562","VIVIANA RODRIGUEZ CORREA","1704",36400,910000),
    @(99,"1143325562","VIVIANA RODRIGUEZ CORREA","1702",36400,910000),
    @(100,"1143325562","VIVIANA RODRIGUEZ CORREA","1701",36400,910000),
    @(101,"1143325562","VIVIANA RODRIGUEZ CORREA","1612",36400,910000),
    @(102,"1143325562","VIVIANA RODRIGUEZ CORREA","1611",36400,910000),
    @(103,"1143325562","VIVIANA RODRIGUEZ CORREA","1610",36400,910000),
    @(104,"80007756","EDGAR ALEXANDER TAUTIVA CASALLAS","1612",42000,1135231),
    @(105,"80007756","EDGAR ALEXANDER TAUTIVA CASALLAS","1611",42000,1135231)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[2]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[3]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[4]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[5]   # G: Salario Basico
}
